# Apply the commit's changes:
#  - Anonymize "fedcore" -> "approach" in the header rows of both sheets
#  - Give the merged-header "spacer" cells (C1/D1 and F1/G1) a light box
#    border (top+bottom, and top+bottom+right for the rightmost spacer)
#  - Drop the stray empty inline-string cell G5 on the computational sheet

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # quality_comparison
$ws2 = $wb.Worksheets.Item(2)   # computational_comparison

# xlLineStyle / Border edge constants used below
$xlContinuous = 1
$xlEdgeLeft   = 7
$xlEdgeTop    = 8
$xlEdgeBottom = 9
$xlEdgeRight  = 10

function Set-SpacerBorders($rng, [bool]$includeRight) {
    # Reset to the plain "Normal" style first so we start from a clean
    # (no font override / no alignment) base, matching the target style.
    $rng.Style = "Normal"
    $rng.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
    $rng.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
    if ($includeRight) {
        $rng.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
    }
}

# --- Sheet 1: quality_comparison ---------------------------------------
Set-SpacerBorders $ws1.Range("C1") $false
Set-SpacerBorders $ws1.Range("D1") $true

$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ----------------------------------
Set-SpacerBorders $ws2.Range("C1") $false
Set-SpacerBorders $ws2.Range("D1") $true
Set-SpacerBorders $ws2.Range("F1") $false
Set-SpacerBorders $ws2.Range("G1") $true

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell G5 entirely.
$ws2.Range("G5").ClearContents()

Write-Output "edit applied"
